$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) and the "SC 92" row (which becomes row 27
# after the first deletion) - matching rows disappear from the diff entirely
# and everything below shifts up.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Update column F (imputed/error values) for the rows that now sit at
# rows 26-33 after the shift.
$ws.Range("F26").Value = 17.38
$ws.Range("F27").ClearContents()
$ws.Range("F30").Value = 16.89
$ws.Range("F32").ClearContents()
